$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "01/04/2013"
$ws.Cells.Item(2, 4).Value = 78.36935839482931

$ws.Cells.Item(3, 3).NumberFormat = "@"
$ws.Cells.Item(3, 3).Value = "01/04/2015"
$ws.Cells.Item(3, 4).Value = 79.11397427123318

$ws.Cells.Item(4, 3).NumberFormat = "@"
$ws.Cells.Item(4, 3).Value = "01/04/2017"
$ws.Cells.Item(4, 4).Value = 79.72812284384679

$ws.Cells.Item(5, 3).NumberFormat = "@"
$ws.Cells.Item(5, 3).Value = "01/04/2019"
$ws.Cells.Item(5, 4).Value = 80.23577424471007

$ws.Cells.Item(6, 3).NumberFormat = "@"
$ws.Cells.Item(6, 3).Value = "01/04/2021"
$ws.Cells.Item(6, 4).Value = 80.73708333931381

$ws.Cells.Item(7, 3).NumberFormat = "@"
$ws.Cells.Item(7, 3).Value = "01/04/2023"
$ws.Cells.Item(7, 4).Value = 81.261256991184

$ws.Cells.Item(8, 3).NumberFormat = "@"
$ws.Cells.Item(8, 3).Value = "01/04/2025"
$ws.Cells.Item(8, 4).Value = 81.83249814315128

$ws.Cells.Item(9, 3).NumberFormat = "@"
$ws.Cells.Item(9, 3).Value = "01/04/2013"
$ws.Cells.Item(9, 4).Value = 45.51490501263478

$ws.Cells.Item(10, 3).NumberFormat = "@"
$ws.Cells.Item(10, 3).Value = "01/04/2015"
$ws.Cells.Item(10, 4).Value = 45.40373533291675

$ws.Cells.Item(11, 3).NumberFormat = "@"
$ws.Cells.Item(11, 3).Value = "01/04/2017"
$ws.Cells.Item(11, 4).Value = 43.63099869346291

$ws.Cells.Item(12, 3).NumberFormat = "@"
$ws.Cells.Item(12, 3).Value = "01/04/2019"
$ws.Cells.Item(12, 4).Value = 44.80581022010754

$ws.Cells.Item(13, 3).NumberFormat = "@"
$ws.Cells.Item(13, 3).Value = "01/04/2021"
$ws.Cells.Item(13, 4).Value = 41.97968547410927

$ws.Cells.Item(14, 3).NumberFormat = "@"
$ws.Cells.Item(14, 3).Value = "01/04/2023"
$ws.Cells.Item(14, 4).Value = 45.99393307422505

$ws.Cells.Item(15, 3).NumberFormat = "@"
$ws.Cells.Item(15, 3).Value = "01/04/2025"
$ws.Cells.Item(15, 4).Value = 48.09755272040089

$ws.Cells.Item(16, 3).NumberFormat = "@"
$ws.Cells.Item(16, 3).Value = "01/04/2013"
$ws.Cells.Item(16, 4).Value = 3.693785299358697

$ws.Cells.Item(17, 3).NumberFormat = "@"
$ws.Cells.Item(17, 3).Value = "01/04/2015"
$ws.Cells.Item(17, 4).Value = 4.169766219339225

$ws.Cells.Item(18, 3).NumberFormat = "@"
$ws.Cells.Item(18, 3).Value = "01/04/2017"
$ws.Cells.Item(18, 4).Value = 6.569385927567932

$ws.Cells.Item(19, 3).NumberFormat = "@"
$ws.Cells.Item(19, 3).Value = "01/04/2019"
$ws.Cells.Item(19, 4).Value = 6.187865072917875

$ws.Cells.Item(20, 3).NumberFormat = "@"
$ws.Cells.Item(20, 3).Value = "01/04/2021"
$ws.Cells.Item(20, 4).Value = 6.971815150253811

$ws.Cells.Item(21, 3).NumberFormat = "@"
$ws.Cells.Item(21, 3).Value = "01/04/2023"
$ws.Cells.Item(21, 4).Value = 4.013650582993648

$ws.Cells.Item(22, 3).NumberFormat = "@"
$ws.Cells.Item(22, 3).Value = "01/04/2025"
$ws.Cells.Item(22, 4).Value = 2.939462030969416

$ws.Cells.Item(23, 3).NumberFormat = "@"
$ws.Cells.Item(23, 3).Value = "01/04/2013"
$ws.Cells.Item(23, 4).Value = 29.16117146049996

$ws.Cells.Item(24, 3).NumberFormat = "@"
$ws.Cells.Item(24, 3).Value = "01/04/2015"
$ws.Cells.Item(24, 4).Value = 29.54096864740481

$ws.Cells.Item(25, 3).NumberFormat = "@"
$ws.Cells.Item(25, 3).Value = "01/04/2017"
$ws.Cells.Item(25, 4).Value = 29.52822756254985

$ws.Cells.Item(26, 3).NumberFormat = "@"
$ws.Cells.Item(26, 3).Value = "01/04/2019"
$ws.Cells.Item(26, 4).Value = 29.24209895168465

$ws.Cells.Item(27, 3).NumberFormat = "@"
$ws.Cells.Item(27, 3).Value = "01/04/2021"
$ws.Cells.Item(27, 4).Value = 31.78558271495075

$ws.Cells.Item(28, 3).NumberFormat = "@"
$ws.Cells.Item(28, 3).Value = "01/04/2023"
$ws.Cells.Item(28, 4).Value = 31.25414731254147

$ws.Cells.Item(29, 3).NumberFormat = "@"
$ws.Cells.Item(29, 3).Value = "01/04/2025"
$ws.Cells.Item(29, 4).Value = 30.79548339178098

$ws.Cells.Item(30, 3).NumberFormat = "@"
$ws.Cells.Item(30, 3).Value = "01/04/2013"
$ws.Cells.Item(30, 4).Value = 49.20818693432935

$ws.Cells.Item(31, 3).NumberFormat = "@"
$ws.Cells.Item(31, 3).Value = "01/04/2015"
$ws.Cells.Item(31, 4).Value = 49.57300562382837

$ws.Cells.Item(32, 3).NumberFormat = "@"
$ws.Cells.Item(32, 3).Value = "01/04/2017"
$ws.Cells.Item(32, 4).Value = 50.20038462103085

$ws.Cells.Item(33, 3).NumberFormat = "@"
$ws.Cells.Item(33, 3).Value = "01/04/2019"
$ws.Cells.Item(33, 4).Value = 50.994158833314

$ws.Cells.Item(34, 3).NumberFormat = "@"
$ws.Cells.Item(34, 3).Value = "01/04/2021"
$ws.Cells.Item(34, 4).Value = 48.95150062436308

$ws.Cells.Item(35, 3).NumberFormat = "@"
$ws.Cells.Item(35, 3).Value = "01/04/2023"
$ws.Cells.Item(35, 4).Value = 50.00710967864252

$ws.Cells.Item(36, 3).NumberFormat = "@"
$ws.Cells.Item(36, 3).Value = "01/04/2025"
$ws.Cells.Item(36, 4).Value = 51.0370147513703

$ws.Cells.Item(37, 3).NumberFormat = "@"
$ws.Cells.Item(37, 3).Value = "01/04/2013"
$ws.Cells.Item(37, 4).Value = 75.96452246726774

$ws.Cells.Item(38, 3).NumberFormat = "@"
$ws.Cells.Item(38, 3).Value = "01/04/2015"
$ws.Cells.Item(38, 4).Value = 76.96066854391862

$ws.Cells.Item(39, 3).NumberFormat = "@"
$ws.Cells.Item(39, 3).Value = "01/04/2017"
$ws.Cells.Item(39, 4).Value = 78.04030951952493

$ws.Cells.Item(40, 3).NumberFormat = "@"
$ws.Cells.Item(40, 3).Value = "01/04/2019"
$ws.Cells.Item(40, 4).Value = 78.6671663900341

$ws.Cells.Item(41, 3).NumberFormat = "@"
$ws.Cells.Item(41, 3).Value = "01/04/2021"

$ws.Cells.Item(42, 3).NumberFormat = "@"
$ws.Cells.Item(42, 3).Value = "01/04/2023"
$ws.Cells.Item(42, 4).Value = 79.97745724802311

$ws.Cells.Item(43, 3).NumberFormat = "@"
$ws.Cells.Item(43, 3).Value = "01/04/2025"
$ws.Cells.Item(43, 4).Value = 80.88755737762359

$ws.Cells.Item(44, 3).NumberFormat = "@"
$ws.Cells.Item(44, 3).Value = "01/04/2013"
$ws.Cells.Item(44, 4).Value = 39.00875920450998

$ws.Cells.Item(45, 3).NumberFormat = "@"
$ws.Cells.Item(45, 3).Value = "01/04/2015"
$ws.Cells.Item(45, 4).Value = 40.19075302025615

$ws.Cells.Item(46, 3).NumberFormat = "@"
$ws.Cells.Item(46, 3).Value = "01/04/2017"
$ws.Cells.Item(46, 4).Value = 36.60968148281447

$ws.Cells.Item(47, 3).NumberFormat = "@"
$ws.Cells.Item(47, 3).Value = "01/04/2019"
$ws.Cells.Item(47, 4).Value = 37.50602345130375

$ws.Cells.Item(48, 3).NumberFormat = "@"
$ws.Cells.Item(48, 3).Value = "01/04/2021"

$ws.Cells.Item(49, 3).NumberFormat = "@"
$ws.Cells.Item(49, 3).Value = "01/04/2023"
$ws.Cells.Item(49, 4).Value = 37.94931403110196

$ws.Cells.Item(50, 3).NumberFormat = "@"
$ws.Cells.Item(50, 3).Value = "01/04/2025"
$ws.Cells.Item(50, 4).Value = 40.44290269455833

$ws.Cells.Item(51, 3).NumberFormat = "@"
$ws.Cells.Item(51, 3).Value = "01/04/2013"
$ws.Cells.Item(51, 4).Value = 4.388783811080302

$ws.Cells.Item(52, 3).NumberFormat = "@"
$ws.Cells.Item(52, 3).Value = "01/04/2015"
$ws.Cells.Item(52, 4).Value = 4.645290217095104

$ws.Cells.Item(53, 3).NumberFormat = "@"
$ws.Cells.Item(53, 3).Value = "01/04/2017"
$ws.Cells.Item(53, 4).Value = 6.947993521684362

$ws.Cells.Item(54, 3).NumberFormat = "@"
$ws.Cells.Item(54, 3).Value = "01/04/2019"
$ws.Cells.Item(54, 4).Value = 6.501757955417537

$ws.Cells.Item(55, 3).NumberFormat = "@"
$ws.Cells.Item(55, 3).Value = "01/04/2021"

$ws.Cells.Item(56, 3).NumberFormat = "@"
$ws.Cells.Item(56, 3).Value = "01/04/2023"
$ws.Cells.Item(56, 4).Value = 4.851975132526726

$ws.Cells.Item(57, 3).NumberFormat = "@"
$ws.Cells.Item(57, 3).Value = "01/04/2025"
$ws.Cells.Item(57, 4).Value = 3.628368197904622

$ws.Cells.Item(58, 3).NumberFormat = "@"
$ws.Cells.Item(58, 3).Value = "01/04/2013"
$ws.Cells.Item(58, 4).Value = 32.56697945167747

$ws.Cells.Item(59, 3).NumberFormat = "@"
$ws.Cells.Item(59, 3).Value = "01/04/2015"
$ws.Cells.Item(59, 4).Value = 32.12462530656735

$ws.Cells.Item(60, 3).NumberFormat = "@"
$ws.Cells.Item(60, 3).Value = "01/04/2017"
$ws.Cells.Item(60, 4).Value = 34.48263451502609

$ws.Cells.Item(61, 3).NumberFormat = "@"
$ws.Cells.Item(61, 3).Value = "01/04/2019"
$ws.Cells.Item(61, 4).Value = 34.66116970962503

$ws.Cells.Item(62, 3).NumberFormat = "@"
$ws.Cells.Item(62, 3).Value = "01/04/2021"

$ws.Cells.Item(63, 3).NumberFormat = "@"
$ws.Cells.Item(63, 3).Value = "01/04/2023"
$ws.Cells.Item(63, 4).Value = 37.17616808439443

$ws.Cells.Item(64, 3).NumberFormat = "@"
$ws.Cells.Item(64, 3).Value = "01/04/2025"
$ws.Cells.Item(64, 4).Value = 36.81628648516066

$ws.Cells.Item(65, 3).NumberFormat = "@"
$ws.Cells.Item(65, 3).Value = "01/04/2013"
$ws.Cells.Item(65, 4).Value = 43.39754301559029

$ws.Cells.Item(66, 3).NumberFormat = "@"
$ws.Cells.Item(66, 3).Value = "01/04/2015"
$ws.Cells.Item(66, 4).Value = 44.83604323735126

$ws.Cells.Item(67, 3).NumberFormat = "@"
$ws.Cells.Item(67, 3).Value = "01/04/2017"
$ws.Cells.Item(67, 4).Value = 43.55767500449883

$ws.Cells.Item(68, 3).NumberFormat = "@"
$ws.Cells.Item(68, 3).Value = "01/04/2019"
$ws.Cells.Item(68, 4).Value = 44.00778140672128

$ws.Cells.Item(69, 3).NumberFormat = "@"
$ws.Cells.Item(69, 3).Value = "01/04/2021"

$ws.Cells.Item(70, 3).NumberFormat = "@"
$ws.Cells.Item(70, 3).Value = "01/04/2023"
$ws.Cells.Item(70, 4).Value = 42.80128916362868

$ws.Cells.Item(71, 3).NumberFormat = "@"
$ws.Cells.Item(71, 3).Value = "01/04/2025"
$ws.Cells.Item(71, 4).Value = 44.07127089246294

$ws.Cells.Item(72, 3).NumberFormat = "@"
$ws.Cells.Item(72, 3).Value = "01/04/2013"
$ws.Cells.Item(72, 4).Value = 76.63246268656717

$ws.Cells.Item(73, 3).NumberFormat = "@"
$ws.Cells.Item(73, 3).Value = "01/04/2015"
$ws.Cells.Item(73, 4).Value = 77.69619091326297

$ws.Cells.Item(74, 3).NumberFormat = "@"
$ws.Cells.Item(74, 3).Value = "01/04/2017"
$ws.Cells.Item(74, 4).Value = 78.46153846153847

$ws.Cells.Item(75, 3).NumberFormat = "@"
$ws.Cells.Item(75, 3).Value = "01/04/2019"
$ws.Cells.Item(75, 4).Value = 78.76620473848904

$ws.Cells.Item(76, 3).NumberFormat = "@"
$ws.Cells.Item(76, 3).Value = "01/04/2021"

$ws.Cells.Item(77, 3).NumberFormat = "@"
$ws.Cells.Item(77, 3).Value = "01/04/2023"
$ws.Cells.Item(77, 4).Value = 78.52760736196319

$ws.Cells.Item(78, 3).NumberFormat = "@"
$ws.Cells.Item(78, 3).Value = "01/04/2025"
$ws.Cells.Item(78, 4).Value = 79.82608695652173

$ws.Cells.Item(79, 3).NumberFormat = "@"
$ws.Cells.Item(79, 3).Value = "01/04/2013"
$ws.Cells.Item(79, 4).Value = 41.18470149253731

$ws.Cells.Item(80, 3).NumberFormat = "@"
$ws.Cells.Item(80, 3).Value = "01/04/2015"
$ws.Cells.Item(80, 4).Value = 42.4506654428637

$ws.Cells.Item(81, 3).NumberFormat = "@"
$ws.Cells.Item(81, 3).Value = "01/04/2017"
$ws.Cells.Item(81, 4).Value = 38.55203619909503

$ws.Cells.Item(82, 3).NumberFormat = "@"
$ws.Cells.Item(82, 3).Value = "01/04/2019"
$ws.Cells.Item(82, 4).Value = 40.32185963343764

$ws.Cells.Item(83, 3).NumberFormat = "@"
$ws.Cells.Item(83, 3).Value = "01/04/2021"

$ws.Cells.Item(84, 3).NumberFormat = "@"
$ws.Cells.Item(84, 3).Value = "01/04/2023"
$ws.Cells.Item(84, 4).Value = 39.83347940403155

$ws.Cells.Item(85, 3).NumberFormat = "@"
$ws.Cells.Item(85, 3).Value = "01/04/2025"
$ws.Cells.Item(85, 4).Value = 41.21739130434783

$ws.Cells.Item(86, 3).NumberFormat = "@"
$ws.Cells.Item(86, 3).Value = "01/04/2013"
$ws.Cells.Item(86, 4).Value = 5.270522388059701

$ws.Cells.Item(87, 3).NumberFormat = "@"
$ws.Cells.Item(87, 3).Value = "01/04/2015"
$ws.Cells.Item(87, 4).Value = 4.313905461220743

$ws.Cells.Item(88, 3).NumberFormat = "@"
$ws.Cells.Item(88, 3).Value = "01/04/2017"
$ws.Cells.Item(88, 4).Value = 6.425339366515836

$ws.Cells.Item(89, 3).NumberFormat = "@"
$ws.Cells.Item(89, 3).Value = "01/04/2019"
$ws.Cells.Item(89, 4).Value = 7.375949932945909

$ws.Cells.Item(90, 3).NumberFormat = "@"
$ws.Cells.Item(90, 3).Value = "01/04/2021"

$ws.Cells.Item(91, 3).NumberFormat = "@"
$ws.Cells.Item(91, 3).Value = "01/04/2023"
$ws.Cells.Item(91, 4).Value = 4.601226993865031

$ws.Cells.Item(92, 3).NumberFormat = "@"
$ws.Cells.Item(92, 3).Value = "01/04/2025"
$ws.Cells.Item(92, 4).Value = 3.608695652173913

$ws.Cells.Item(93, 3).NumberFormat = "@"
$ws.Cells.Item(93, 3).Value = "01/04/2013"
$ws.Cells.Item(93, 4).Value = 30.17723880597015

$ws.Cells.Item(94, 3).NumberFormat = "@"
$ws.Cells.Item(94, 3).Value = "01/04/2015"
$ws.Cells.Item(94, 4).Value = 30.9775126204681

$ws.Cells.Item(95, 3).NumberFormat = "@"
$ws.Cells.Item(95, 3).Value = "01/04/2017"
$ws.Cells.Item(95, 4).Value = 33.43891402714932

$ws.Cells.Item(96, 3).NumberFormat = "@"
$ws.Cells.Item(96, 3).Value = "01/04/2019"
$ws.Cells.Item(96, 4).Value = 31.0683951721055

$ws.Cells.Item(97, 3).NumberFormat = "@"
$ws.Cells.Item(97, 3).Value = "01/04/2021"

$ws.Cells.Item(98, 3).NumberFormat = "@"
$ws.Cells.Item(98, 3).Value = "01/04/2023"
$ws.Cells.Item(98, 4).Value = 34.09290096406661

$ws.Cells.Item(99, 3).NumberFormat = "@"
$ws.Cells.Item(99, 3).Value = "01/04/2025"
$ws.Cells.Item(99, 4).Value = 35.04347826086956

$ws.Cells.Item(100, 3).NumberFormat = "@"
$ws.Cells.Item(100, 3).Value = "01/04/2013"
$ws.Cells.Item(100, 4).Value = 46.50186567164179

$ws.Cells.Item(101, 3).NumberFormat = "@"
$ws.Cells.Item(101, 3).Value = "01/04/2015"
$ws.Cells.Item(101, 4).Value = 46.76457090408444

$ws.Cells.Item(102, 3).NumberFormat = "@"
$ws.Cells.Item(102, 3).Value = "01/04/2017"
$ws.Cells.Item(102, 4).Value = 45.02262443438914

$ws.Cells.Item(103, 3).NumberFormat = "@"
$ws.Cells.Item(103, 3).Value = "01/04/2019"
$ws.Cells.Item(103, 4).Value = 47.74251229324989

$ws.Cells.Item(104, 3).NumberFormat = "@"
$ws.Cells.Item(104, 3).Value = "01/04/2021"

$ws.Cells.Item(105, 3).NumberFormat = "@"
$ws.Cells.Item(105, 3).Value = "01/04/2023"
$ws.Cells.Item(105, 4).Value = 44.3908851884312

$ws.Cells.Item(106, 3).NumberFormat = "@"
$ws.Cells.Item(106, 3).Value = "01/04/2025"
$ws.Cells.Item(106, 4).Value = 44.82608695652173
